$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.45 = 9252.59 pesos`n✅ 9252.59 pesos = 2.44 = 943.56 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 407.995
$ws2.Range("O10").Value = 3775.01
$ws2.Range("N12").Value = 3795
$ws2.Range("O12").Value = 387.005
